# Updated simulation files with Holden scheme
#
# The sheet originally had the 18 HKL/Pairs columns (C:T) repeated three
# times (C:T, U:AD, AE:AN) for three duplicate "groups", plus rows for the
# BT8Hex / Spiral / Offset / HexGrid scan schemes. This edit:
#   1) drops the two duplicated column blocks (U:AN), keeping one copy
#   2) reorders the HKL-index header labels (C2:J2)
#   3) inserts four new "Holden" scheme rows (2.5/5/10/15 deg) right
#      before the HexGrid rows, with the same all-1s data as every other
#      scheme row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the two duplicated column blocks (U:AN), leaving a single
#    18-column data block (C:T).
$ws.Range("U1:AN19").EntireColumn.Delete()

# 2) Re-order the HKL bracket-index headers in row 2 to match the new
#    layout.
$hklHeaders = @("[4, 0, 0]", "[2, 1, 1]", "[2, 2, 0]", "[2, 0, 0]", "[2, 2, 2]", "[3, 1, 0]", "[1, 1, 0]", "[3, 2, 1]")
for ($i = 0; $i -lt $hklHeaders.Length; $i++) {
    $ws.Cells.Item(2, 3 + $i).Value = $hklHeaders[$i]
}

# 3) Insert 4 new rows right above the first HexGrid row (row 16) to hold
#    the new Holden scheme entries, pushing the HexGrid rows down to 20:23.
$ws.Rows("16:19").Insert()

# 4) Populate the new Holden rows (column A index values are fixed up by
#    the renumbering pass below).
$holdenNames = @("Holden2.5", "Holden5", "Holden10", "Holden15")
for ($i = 0; $i -lt $holdenNames.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 2).Value = $holdenNames[$i]
    for ($c = 3; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

# Match the bold / bordered / centered style used by the rest of column A
# (and carry it onto the new rows' A cells).
$ws.Range("A15").Copy()
$ws.Range("A16:A19").PasteSpecial(-4122)

# Renumber column A (the running HKL index) for every data row so it stays
# a contiguous 0-based sequence after the insert.
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

$excel.CutCopyMode = 0
